{"js": "// 1) Post-conditions: \"A location/region is defined on the map.\" ->\n//    \"A location is defined on the map.\"  (drop \"/region\")\nconst postCond = context.document.body.search(\n  \"A location/region is defined on the map.\",\n  { matchCase: true }\n);\npostCond.load(\"items\");\nawait context.sync();\nif (postCond.items.length > 0) {\n  postCond.items[0].insertText(\"A location is defined on the map.\", \"Replace\");\n}\n\n// 2) Main path step 1: \"... Search Map field ...\" ->\n//    \"... Search the map field ...\"\nconst searchStep = context.document.body.search(\n  \"enters the name of a location in the Search Map field and presses the search button\",\n  { matchCase: true }\n);\nsearchStep.load(\"items\");\nawait context.sync();\nif (searchStep.items.length > 0) {\n  searchStep.items[0].insertText(\n    \"enters the name of a location in the Search the map field and presses the search button\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n\n// 3) Remove the whole \"Alternate paths / R26-2 ...\" table row (it was\n//    dropped entirely from the use case).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of table.rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nlet targetRow = null;\nfor (const row of table.rows.items) {\n  const rowText = row.cells.items.map((c) => c.body.text).join(\" \");\n  if (rowText.indexOf(\"Alternate paths\") !== -1 && rowText.indexOf(\"R26-2\") !== -1) {\n    targetRow = row;\n    break;\n  }\n}\n\nif (targetRow) {\n  targetRow.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Post-conditions: \"A location/region is defined on the map.\" ->\n#    \"A location is defined on the map.\"  (drop \"/region\")\n$find1 = $d.Content\n$find1.Find.Text = \"A location/region is defined on the map.\"\n$find1.Find.MatchCase = $true\nif ($find1.Find.Execute()) {\n    $find1.Text = \"A location is defined on the map.\"\n}\n\n# 2) Main path step 1: \"... Search Map field ...\" ->\n#    \"... Search the map field ...\"\n$find2 = $d.Content\n$find2.Find.Text = \"enters the name of a location in the Search Map field and presses the search button\"\n$find2.Find.MatchCase = $true\nif ($find2.Find.Execute()) {\n    $find2.Text = \"enters the name of a location in the Search the map field and presses the search button\"\n}\n\n# 3) Remove the whole \"Alternate paths / R26-2 ...\" table row (it was\n#    dropped entirely from the use case).\n$table = $d.Tables.Item(1)\nfor ($i = $table.Rows.Count; $i -ge 1; $i--) {\n    $row = $table.Rows.Item($i)\n    $rowText = $row.Range.Text\n    if ($rowText -like \"*Alternate paths*\" -and $rowText -like \"*R26-2*\") {\n        $row.Delete()\n        break\n    }\n}\n"}
